$d = $word.ActiveDocument

# 1) Header/title run: standalone "2020" -> "2021"
$d.Content.Find.Execute("2020", $true, $false, $false, $false, $false, $true, 1, $false, "2021", 2)

# 2) Month name: "نوفمبر" (November) -> "جانفي" (January)
$d.Content.Find.Execute("نوفمبر", $true, $false, $false, $false, $false, $true, 1, $false, "جانفي", 2)

# 3) Amount: "41 770 000,00" -> "44 260 000,00"
$d.Content.Find.Execute("41 770 000,00", $true, $false, $false, $false, $false, $true, 1, $false, "44 260 000,00", 2)
